$d = $word.ActiveDocument

# 1. Remove the (now orphaned) heading bookmark anchored at the start of the
#    "AR ${AR}" heading run. It is a hidden "_heading=h.*" bookmark that
#    Word does not surface in Bookmarks.Count, but it can still be located
#    and removed by name.
try {
    $d.Bookmarks.Item("_heading=h.gjdgxs").Delete()
} catch {
}

# 2. Heading paragraph: "AR ${AR}" -> "Recommendation ${REC}"
#    Two scoped, whole-word-safe replacements so the "${AR}" placeholder and
#    the bare "AR" acronym are handled distinctly (and "ARC Number" later in
#    the document is left untouched).
$headingRange = $d.Paragraphs.Item(1).Range
$headingRange.Find.Execute("AR", $true, $true, $false, $false, $false, $true, 1, $false, "Recommendation", 1)

$headingRange = $d.Paragraphs.Item(1).Range
$headingRange.Find.Execute("`$`{AR`}", $true, $false, $false, $false, $false, $true, 1, $false, "`$`{REC`}", 1)

# 3. Body paragraph: "...implementation of this AR, any effects..."
#    -> "...implementation of this recommendation, any effects..."
$d.Content.Find.Execute("this AR,", $true, $false, $false, $false, $false, $true, 1, $false, "this recommendation,", 2)

# 4. Body paragraph: "...savings for this AR is " -> "...savings for this recommendation is "
$d.Content.Find.Execute("this AR is ", $true, $false, $false, $false, $false, $true, 1, $false, "this recommendation is ", 2)
